$wb = $excel.ActiveWorkbook

$successMsg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet "Bico": rename Obs -> Obs_relatorio, add new Obs_sped column ---
$wsBico = $wb.Worksheets.Item("Bico")

$wsBico.Range("H1").Value = "Obs_relatorio"
$wsBico.Range("I1").Value = "Obs_sped"

for ($row = 2; $row -le 9; $row++) {
    $wsBico.Cells.Item($row, 8).Value = $successMsg
}

# --- Sheet "Tanque": rename Obs -> Obs_relatorio, add new Obs_sped column ---
$wsTanque = $wb.Worksheets.Item("Tanque")

$wsTanque.Range("F1").Value = "Obs_relatorio"
$wsTanque.Range("G1").Value = "Obs_sped"

for ($row = 2; $row -le 5; $row++) {
    $wsTanque.Cells.Item($row, 6).Value = $successMsg
}
